$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 313, shifting existing rows 313:332 down to 314:333
$ws.Range("A313").EntireRow.Insert()

# Populate the newly inserted row 313 with the new weekly price observation
$ws.Cells.Item(313, 1).Value = 9
$ws.Cells.Item(313, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(313, 3).Value = "Metropolitana"
$ws.Cells.Item(313, 4).Value = 44578
$ws.Cells.Item(313, 5).Value = 13
$ws.Cells.Item(313, 6).Value = "Fruta"
$ws.Cells.Item(313, 7).Value = 100108
$ws.Cells.Item(313, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(313, 9).Value = 100108002
$ws.Cells.Item(313, 10).Value = "Mango"
$ws.Cells.Item(313, 11).Value = "Sin especificar"
$ws.Cells.Item(313, 12).Value = "Primera"
$ws.Cells.Item(313, 13).Value = 590
$ws.Cells.Item(313, 14).Value = 5000
$ws.Cells.Item(313, 15).Value = 6000
$ws.Cells.Item(313, 16).Value = 5644
$ws.Cells.Item(313, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(313, 18).Value = "Perú"
$ws.Cells.Item(313, 19).Value = 1411
$ws.Cells.Item(313, 20).Value = 4

# Apply the same date-cell number format ("D" column) the rest of the sheet uses
$ws.Cells.Item(313, 4).NumberFormat = $ws.Cells.Item(314, 4).NumberFormat
